$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.013315305999535225
$ws.Range("B3").Value = 0.073819035627262664
$ws.Range("B4").Value = 0.029670863717644842
$ws.Range("B5").Value = 0.018931464752434613
$ws.Range("B6").Value = 0.26403725212990764
$ws.Range("B7").Value = 0.085392997239621141
$ws.Range("B8").Value = 0.023448441112473167
$ws.Range("B9").Value = 0.3214577998636674
